# Update the work-report worksheet so it reflects the single-work-request
# reconciled numbers (pricing populated, duplicate/incorrect Monday line
# items removed, totals + header summary refreshed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header / summary block -------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:00 AM"
$ws.Range("C8").Value = 5160.76
$ws.Range("C9").Value = 19
$ws.Range("G10").Value = ""

# ---- Monday (07/28/2025) line items -----------------------------------------
# The old rows 16-37 had 22 line items (several bogus/duplicate ones, all
# priced at 0). The corrected sheet has 16 line items with real pricing.
# Easiest correct approach: clear the old block (16 down through the old
# second TOTAL row 46) then write the new, shorter block back in.

$ws.Range("A16:I46").Clear()

$monday = @(
  @("Point 01","SWI-15-CO1-100-H","Inst","SWI,15kV,Line Cutout 1PH,100A,Hook","EA",3,187.26),
  @("Point 02","SWI-15-CO1-100-H","Rem","SWI,15kV,Line Cutout 1PH,100A,Hook","EA",3,79.56),
  @("Point 07","ARM-8SF-GN-TL","Inst","ARM,8ft Sgl.Fiberglass,Gain,Tangent LD","EA",1,350.53),
  @("Point 07","INS-15-P-S","Inst","INS,15kV,Pin,Silicon Polymer","EA",3,282.51),
  @("Point 07","PIN-35-PTP","Inst","Pin,35kV,Pole Top","EA",1,94.17),
  @("Point 07","PIN-XAL","Inst","Pin,Crossarm Light","EA",2,188.34),
  @("Point 07","SAA-3-CV","Inst","SAA,3 inch,Clevis","EA",1,55.18),
  @("Point 08","ARM-8SF-GN-TL","Rem","ARM,8ft Sgl.Fiberglass,Gain,Tangent LD","EA",1,61.83),
  @("Point 08","INS-15-P-S","Rem","INS,15kV,Pin,Silicon Polymer","EA",3,93.23999999999999),
  @("Point 08","PIN-35-PTP","Rem","Pin,35kV,Pole Top","EA",1,31.08),
  @("Point 08","PIN-XAL","Rem","Pin,Crossarm Light","EA",2,62.16),
  @("Point 08","PLA-CUT","Rem","PLA,Cut Off Pole Top","EA",1,216.17),
  @("Point 08","POL-40-2","Rem","Pole,40ft,Class 2","EA",1,198.88),
  @("Point 08","SAA-3-CV","Rem","SAA,3 inch,Clevis","EA",1,17.2),
  @("Point 03","PLA-HDIG","Inst","PLA,Hand Dig or Additional  Excavation","EA",1,648.53),
  @("Point 05","PLA-HDIG","Inst","PLA,Hand Dig or Additional  Excavation","EA",1,648.53)
)

$row = 16
foreach ($item in $monday) {
  $ws.Cells.Item($row, 1).Value = $item[0]
  $ws.Cells.Item($row, 2).Value = $item[1]
  $ws.Cells.Item($row, 3).Value = $item[2]
  $ws.Cells.Item($row, 4).Value = $item[3]
  $ws.Cells.Item($row, 5).Value = $item[4]
  $ws.Cells.Item($row, 6).Value = $item[5]
  $ws.Cells.Item($row, 8).Value = $item[6]
  $row++
}

$mondayTotalRow = $row
$ws.Cells.Item($mondayTotalRow, 1).Value = "TOTAL"
$ws.Cells.Item($mondayTotalRow, 8).Value = 3215.17

# ---- Tuesday (07/29/2025) section -------------------------------------------
$tueHeaderRow = $mondayTotalRow + 3
$ws.Cells.Item($tueHeaderRow, 1).Value = "Tuesday (07/29/2025)"

$tueColRow = $tueHeaderRow + 1
$ws.Cells.Item($tueColRow, 1).Value = "Point Number"
$ws.Cells.Item($tueColRow, 2).Value = "Billable Unit Code"
$ws.Cells.Item($tueColRow, 3).Value = "Work Type"
$ws.Cells.Item($tueColRow, 4).Value = "Unit Description"
$ws.Cells.Item($tueColRow, 5).Value = "Unit of Measure"
$ws.Cells.Item($tueColRow, 6).Value = "# Units"
$ws.Cells.Item($tueColRow, 7).Value = "N/A"
$ws.Cells.Item($tueColRow, 8).Value = "Pricing"

$tuesday = @(
  @("Point 05","PLA-HDIG","Inst","PLA,Hand Dig or Additional  Excavation","EA",1,648.53),
  @("Point 07","PLA-HDIG","Inst","PLA,Hand Dig or Additional  Excavation","EA",1,648.53),
  @("Point 09","PLA-HDIG","Inst","PLA,Hand Dig or Additional  Excavation","EA",1,648.53)
)

$row = $tueColRow + 1
foreach ($item in $tuesday) {
  $ws.Cells.Item($row, 1).Value = $item[0]
  $ws.Cells.Item($row, 2).Value = $item[1]
  $ws.Cells.Item($row, 3).Value = $item[2]
  $ws.Cells.Item($row, 4).Value = $item[3]
  $ws.Cells.Item($row, 5).Value = $item[4]
  $ws.Cells.Item($row, 6).Value = $item[5]
  $ws.Cells.Item($row, 8).Value = $item[6]
  $row++
}

$tueTotalRow = $row
$ws.Cells.Item($tueTotalRow, 1).Value = "TOTAL"
$ws.Cells.Item($tueTotalRow, 8).Value = 1945.59
